$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H3").Value = 6
$ws.Range("I3").Value = 7.6
$ws.Range("V3").Value = 1.15
$ws.Range("Y3").Value = 28
$ws.Range("AB3").Value = 970
$ws.Range("AD3").Value = 34
$ws.Range("AF3").Value = 970
$ws.Range("AG3").Value = 12
$ws.Range("AH3").Value = 28
$ws.Range("AJ3").Value = 970
$ws.Range("AN3").Value = 10
$ws.Range("F4").Value = 5.6
$ws.Range("M4").Value = 1.04
$ws.Range("N4").Value = 5.3
$ws.Range("P4").Value = 2.42
$ws.Range("S4").Value = 2.68
$ws.Range("X4").Value = 23
$ws.Range("AJ4").Value = 150
$ws.Range("AO4").Value = 7.2
$ws.Range("G5").Value = 2.34
$ws.Range("H5").Value = 3.3
$ws.Range("H6").Value = 4.5
$ws.Range("O6").Value = 1.16
$ws.Range("S6").Value = 2.1
$ws.Range("T6").Value = 1.56
$ws.Range("Z6").Value = 46
$ws.Range("AB6").Value = 14.5
$ws.Range("AD6").Value = 21
$ws.Range("AE6").Value = 55
$ws.Range("AF6").Value = 14.5
$ws.Range("AG6").Value = 11.5
$ws.Range("AM6").Value = 65
$ws.Range("AN6").Value = 6.8
$ws.Range("AO6").Value = 38
$ws.Range("R7").Value = 2.02
$ws.Range("T7").Value = 1.65
$ws.Range("U7").Value = 2.48
$ws.Range("X7").Value = 42
$ws.Range("AF7").Value = 95
$ws.Range("F8").Value = 2.94
$ws.Range("H8").Value = 2.48
$ws.Range("AC8").Value = 7.8
$ws.Range("S9").Value = 1.32
$ws.Range("F10").Value = 5.1
$ws.Range("I10").Value = 1.66
$ws.Range("K10").Value = 5.6
$ws.Range("M10").Value = 1.06
$ws.Range("N10").Value = 1.01
$ws.Range("V10").Value = 2.5
$ws.Range("W10").Value = 1.14
$ws.Range("Y10").Value = 970
$ws.Range("AB10").Value = 25
$ws.Range("AE10").Value = 970
$ws.Range("AG10").Value = 30
$ws.Range("AH10").Value = 29
$ws.Range("L11").Value = 1.2
$ws.Range("N11").Value = 5.3
$ws.Range("P11").Value = 2.46
$ws.Range("Q11").Value = 1.54
$ws.Range("T11").Value = 1.01
$ws.Range("F13").Value = 3.15
$ws.Range("I13").Value = 2.36
$ws.Range("J13").Value = 3.65
$ws.Range("K13").Value = 4.2
$ws.Range("P13").Value = 2.08
$ws.Range("Q13").Value = 1.75
$ws.Range("S13").Value = 2.86
$ws.Range("T13").Value = 1.65
$ws.Range("U13").Value = 2.26
$ws.Range("V13").Value = 1.73
$ws.Range("AA13").Value = 32
$ws.Range("AC13").Value = 9.199999999999999
$ws.Range("AF13").Value = 26
$ws.Range("AL13").Value = 44
$ws.Range("AM13").Value = 80
$ws.Range("AN13").Value = 30
$ws.Range("I14").Value = 13.5
$ws.Range("U14").Value = 1.64
$ws.Range("V14").Value = 1.08
$ws.Range("X14").Value = 29
$ws.Range("AB14").Value = 12.5
$ws.Range("AF14").Value = 12.5
$ws.Range("G16").Value = 2.82
$ws.Range("I16").Value = 2.62
$ws.Range("J16").Value = 3.75
$ws.Range("N16").Value = 4.5
$ws.Range("U16").Value = 2.34
$ws.Range("AA16").Value = 38
$ws.Range("AJ16").Value = 42
$ws.Range("H17").Value = 1.91
$ws.Range("P17").Value = 1.85
$ws.Range("T17").Value = 1.93
$ws.Range("W17").Value = 1.26
$ws.Range("AB17").Value = 15
$ws.Range("F18").Value = 2.96
$ws.Range("G18").Value = 3
$ws.Range("H18").Value = 2.6
$ws.Range("I18").Value = 2.62
$ws.Range("L18").Value = 1.46
$ws.Range("N18").Value = 3.7
$ws.Range("O18").Value = 1.35
$ws.Range("P18").Value = 1.89
$ws.Range("Q18").Value = 2.08
$ws.Range("R18").Value = 1.34
$ws.Range("S18").Value = 3.75
$ws.Range("T18").Value = 1.81
$ws.Range("V18").Value = 1.61
$ws.Range("W18").Value = 1.5
$ws.Range("Y18").Value = 10.5
$ws.Range("Z18").Value = 16.5
$ws.Range("AI18").Value = 44
$ws.Range("AJ18").Value = 48
$ws.Range("AK18").Value = 34
$ws.Range("AL18").Value = 46
$ws.Range("AN18").Value = 34
$ws.Range("G19").Value = 4.7
$ws.Range("H19").Value = 1.8
$ws.Range("I19").Value = 1.82
$ws.Range("N19").Value = 5.4
$ws.Range("P19").Value = 2.48
$ws.Range("AF19").Value = 38
$ws.Range("AI19").Value = 27
$ws.Range("AN19").Value = 40
$ws.Range("F20").Value = 1.56
$ws.Range("G20").Value = 1.58
$ws.Range("H20").Value = 6
$ws.Range("I20").Value = 6.2
$ws.Range("K20").Value = 5.1
$ws.Range("S20").Value = 2.18
$ws.Range("T20").Value = 1.6
$ws.Range("U20").Value = 2.58
$ws.Range("V20").Value = 1.19
$ws.Range("W20").Value = 2.72
$ws.Range("AB20").Value = 14
$ws.Range("AD20").Value = 23
$ws.Range("AE20").Value = 65
$ws.Range("AH20").Value = 17.5
$ws.Range("AI20").Value = 55
$ws.Range("AJ20").Value = 16
$ws.Range("AN20").Value = 5.3
$ws.Range("AO20").Value = 46
$ws.Range("F21").Value = 1.3
$ws.Range("G21").Value = 1.32
$ws.Range("I21").Value = 11.5
$ws.Range("K21").Value = 7
$ws.Range("P21").Value = 3.4
$ws.Range("R21").Value = 1.97
$ws.Range("S21").Value = 1.99
$ws.Range("T21").Value = 1.73
$ws.Range("U21").Value = 2.24
$ws.Range("W21").Value = 4.2
$ws.Range("Z21").Value = 120
$ws.Range("AB21").Value = 14.5
$ws.Range("AE21").Value = 130
$ws.Range("AI21").Value = 990
$ws.Range("AN21").Value = 3.65
$ws.Range("M22").Value = 1.06
$ws.Range("AO22").Value = 20
